$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("C7").Value = 8000063

# Rows 20-25
$ws.Range("C20").Value = 8000035
$ws.Range("C21").Value = 8000038
$ws.Range("C22").Value = 8000039
$ws.Range("C23").Value = 8000042
$ws.Range("C24").Value = 8000052
$ws.Range("C25").Value = 8000058

# Rows 28-29
$ws.Range("C28").Value = 8000065
$ws.Range("C29").Value = 8000070

# Row 30
$ws.Range("B30").Value = 10000036
$ws.Range("C30").Value = 8000036
$ws.Range("D30").Value = "Repubblica e Cantone Ticino"
$ws.Range("E30").Value = "Legge sulla protezione civile del 26 febbraio 2007"

# Row 31
$ws.Range("B31").Value = 10000036
$ws.Range("C31").Value = 8000047
$ws.Range("D31").Value = "Repubblica e Cantone Ticino"
$ws.Range("E31").Value = "Legge sulla protezione della popolazione (del 26 febbraio 2007)"

# Row 32
$ws.Range("B32").Value = 10000036
$ws.Range("C32").Value = 8000054
$ws.Range("D32").Value = "Repubblica e Cantone Ticino"
$ws.Range("E32").Value = "Regolamento sulla protezione della popolazione (RProtPop) (del 18 ottobre 2017)"

# Row 33
$ws.Range("B33").Value = 10000036
$ws.Range("C33").Value = 8000062
$ws.Range("E33").Value = "Servizio della protezione della popolazione"

# Row 34
$ws.Range("C34").Value = 8000041
$ws.Range("D34").Value = "IRPI CNR"
$ws.Range("E34").Value = "Modelli e carte di suscettibilità da frana"

# Row 35
$ws.Range("B35").Value = 10000045
$ws.Range("C35").Value = 8000042
$ws.Range("D35").Value = "Confederazione elvetica"
$ws.Range("E35").Value = "Legge federale sulla protezione della popolazione e sulla protezione civile del 4 ottobre 2002"

# Row 36
$ws.Range("B36").Value = 10000073
$ws.Range("C36").Value = 8000073
$ws.Range("D36").Value = "Areu Lombardia"
$ws.Range("E36").Value = "Le SOREU"

# Row 37
$ws.Range("B37").Value = 10000073
$ws.Range("C37").Value = 8000074
$ws.Range("D37").Value = "Areu Lombardia"
$ws.Range("E37").Value = "SOREU dei Laghi"

# Row 38
$ws.Range("B38").Value = 10000075
$ws.Range("C38").Value = 8000075
$ws.Range("D38").Value = "Repubblica e Canton Ticino"
$ws.Range("E38").Value = "Chi siamo"
